$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.335.52'
$ws.Range('D3').Value = '1.567.37'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.06%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.97'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  -0.78%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.42'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -3.97%  '
$ws.Range('E9').Value = '  -1.84%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.245'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '1.790.43'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = '1.571.52'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '28.324.14'
$ws.Range('E16').Value = '  -0.78%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.514'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -1.19%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.12'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -1.42%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.71'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.18%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.39'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').Value = '0.0₃0678'
$ws.Range('E21').Value = '  -2.16%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  +1.90%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.94'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -2.08%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.05'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('D35').Value = '1.381.48'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('E42').Value = '  +3.74%  '
$ws.Range('E43').Value = '  +0.00%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0476'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -0.35%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.784'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('E46').Value = '  -3.36%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.25'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.60%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.918'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -6.17%  '
$ws.Range('D49').Value = '1.703.46'
$ws.Range('E49').Value = '  +0.15%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.60'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').Value = '0.0₆0100'
$ws.Range('E51').Value = '  -2.40%  '
